# Recompute of the DOMLEM rule-induction run produced a fresh set of
# results (rule object lists came back from a Python set() in a new,
# but equivalent, order; cross-validation + per-rule statistics were
# refreshed as well).

$wb = $excel.ActiveWorkbook

# --- "Reguły" sheet (8th sheet) : rule descriptions in column B ---------
$wsRules = $wb.Worksheets.Item(8)

$wsRules.Range("B2").Value = "(age >=  40.0) & (woman_eval <=  2.0) => (class <= 1) ['a49', 'a23']"
$wsRules.Range("B3").Value = "(age >=  38.0) => (class <= 2) ['a17', 'a48', 'a11', 'a23', 'a6', 'a41', 'a33', 'a36', 'a3', 'a5', 'a51', 'a49']"
$wsRules.Range("B4").Value = "(morpho_quality <=  9.0) => (class <= 2) ['a30', 'a7', 'a11', 'a9', 'a33', 'a4']"
$wsRules.Range("B5").Value = "(age >=  27.0) => (class <= 3) ['a26', 'a9', 'a18', 'a33', 'a38', 'a34', 'a14', 'a51', 'a29', 'a49', 'a12', 'a30', 'a40', 'a17', 'a24', 'a23', 'a25', 'a32', 'a43', 'a5', 'a21', 'a4', 'a44', 'a2', 'a10', 'a8', 'a20', 'a47', 'a48', 'a6', 'a37', 'a36', 'a16', 'a31', 'a35', 'a50', 'a7', 'a11', 'a39', 'a46', 'a45', 'a28', 'a41', 'a1', 'a27', 'a3', 'a15', 'a13', 'a19', 'a22']"
$wsRules.Range("B6").Value = "(sperm <=  2.0) => (class <= 3) ['a5', 'a42', 'a2', 'a7']"
$wsRules.Range("B7").Value = "(age <=  38.0) => (class >= 3) ['a26', 'a9', 'a18', 'a38', 'a34', 'a14', 'a29', 'a12', 'a30', 'a40', 'a17', 'a24', 'a25', 'a32', 'a43', 'a21', 'a4', 'a44', 'a2', 'a10', 'a8', 'a20', 'a47', 'a37', 'a36', 'a16', 'a31', 'a35', 'a50', 'a7', 'a39', 'a46', 'a42', 'a45', 'a28', 'a41', 'a1', 'a27', 'a3', 'a15', 'a13', 'a19', 'a22']"
$wsRules.Range("B8").Value = "(infertility <=  2.0) => (class >= 3) ['a18', 'a38', 'a14', 'a49', 'a12', 'a25', 'a43', 'a2', 'a10', 'a8', 'a47', 'a48', 'a7', 'a39', 'a46', 'a45', 'a1', 'a3', 'a13', 'a22']"
$wsRules.Range("B9").Value = "(age <=  39.0) => (class >= 2) ['a26', 'a9', 'a18', 'a33', 'a38', 'a34', 'a14', 'a29', 'a12', 'a30', 'a40', 'a17', 'a24', 'a25', 'a32', 'a43', 'a21', 'a4', 'a44', 'a2', 'a10', 'a8', 'a20', 'a47', 'a37', 'a36', 'a16', 'a31', 'a35', 'a50', 'a7', 'a11', 'a39', 'a46', 'a42', 'a45', 'a28', 'a41', 'a1', 'a27', 'a3', 'a15', 'a13', 'a19', 'a22']"
$wsRules.Range("B10").Value = "(woman_eval >=  3.0) => (class >= 2) ['a26', 'a18', 'a38', 'a34', 'a14', 'a51', 'a29', 'a12', 'a30', 'a40', 'a17', 'a24', 'a25', 'a32', 'a43', 'a5', 'a21', 'a4', 'a44', 'a2', 'a10', 'a8', 'a20', 'a47', 'a48', 'a6', 'a37', 'a36', 'a16', 'a50', 'a7', 'a11', 'a39', 'a45', 'a41', 'a1', 'a27', 'a13', 'a19']"
$wsRules.Range("B11").Value = "(infertility <=  2.0) => (class >= 2) ['a18', 'a38', 'a14', 'a49', 'a12', 'a25', 'a43', 'a2', 'a10', 'a8', 'a47', 'a48', 'a7', 'a39', 'a46', 'a45', 'a1', 'a3', 'a13', 'a22']"

# --- "Statystyki reguł" sheet (9th sheet) : refreshed coverage column C -
$wsStats = $wb.Worksheets.Item(9)

$wsStats.Range("C2").Value = 2
$wsStats.Range("C4").Value = 0.6666666666666666
$wsStats.Range("C6").Value = 0.08
$wsStats.Range("C10").Value = 0.78

# --- "Walidacja krzyżowa" sheet (10th sheet) : metric rows re-labelled -
# (the underlying shared-string table was reshuffled so the same four
# metric rows now read accuracy / not_classified / correct / f1_score)
$wsCV = $wb.Worksheets.Item(10)

$wsCV.Range("A1").Value = "accuracy"
$wsCV.Range("B1").Value = 0.7450980392156863
$wsCV.Range("A2").Value = "not_classified"
$wsCV.Range("B2").Value = 0
$wsCV.Range("A3").Value = "correct"
$wsCV.Range("B3").Value = 0.7450980392156863
$wsCV.Range("A4").Value = "f1_score"
$wsCV.Range("B4").Value = 0
